$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before row 570, shifting existing rows 570-661 down to 574-665.
$ws.Range("A570:A573").EntireRow.Insert()

# Constant columns shared by every data row in this table.
$marketId = 6
$marketName = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$category = 13
$fruitLabel = "Fruta"
$groupId = 100108
$groupName = "Tropicales y subtropicales"
$productId = 100108006
$productName = "Plátano"
$unit = "$/caja 20 kilos"
$origin = "Ecuador"
$kgPerUnit = 20

$newRows = @(
    @{ Row = 570; Date = 44474; Variety = "Sin especificar"; Quality = "Pintón";         Vol = 800;  Min = 17000; Max = 18000; Avg = 17500; PriceKg = 875 },
    @{ Row = 571; Date = 44474; Variety = "Sin especificar"; Quality = "Primera Maduro"; Vol = 380;  Min = 20000; Max = 20000; Avg = 20000; PriceKg = 1000 },
    @{ Row = 572; Date = 44474; Variety = "Sin especificar"; Quality = "Primera Pintón"; Vol = 2160; Min = 18000; Max = 20000; Avg = 19074; PriceKg = 954 },
    @{ Row = 573; Date = 44474; Variety = "Sin especificar"; Quality = "Primera Verde";  Vol = 1670; Min = 18000; Max = 20000; Avg = 19138; PriceKg = 957 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $marketId
    $ws.Cells.Item($row, 2).Value = $marketName
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.Date
    $ws.Cells.Item($row, 5).Value = $category
    $ws.Cells.Item($row, 6).Value = $fruitLabel
    $ws.Cells.Item($row, 7).Value = $groupId
    $ws.Cells.Item($row, 8).Value = $groupName
    $ws.Cells.Item($row, 9).Value = $productId
    $ws.Cells.Item($row, 10).Value = $productName
    $ws.Cells.Item($row, 11).Value = $r.Variety
    $ws.Cells.Item($row, 12).Value = $r.Quality
    $ws.Cells.Item($row, 13).Value = $r.Vol
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Avg
    $ws.Cells.Item($row, 17).Value = $unit
    $ws.Cells.Item($row, 18).Value = $origin
    $ws.Cells.Item($row, 19).Value = $r.PriceKg
    $ws.Cells.Item($row, 20).Value = $kgPerUnit
}
